$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "CellRef=NewValue" pairs taken from the 02.08.2021 daily-stats refresh:
#  - backfilled AgTests/AgPosit (F/G) for 2020-06-19 (row 118)
#  - corrected several historical AgTests/AgPosit figures (rows 264-511)
#  - appended four new daily rows for 2021-07-29 .. 2021-08-01 (rows 512-515)
$updates = @(
"F118=12","G118=0","F264=43113","G264=884","F271=45543","G271=1717","F272=30271","G272=1597",
"F273=32344","G273=1675","F274=28876","G274=1274","F275=30658","G275=1260","F276=10764","G276=367",
"F278=31880","G278=2164","F279=42544","G279=2960","F280=34946","G280=2311","F281=47097","G281=3186",
"F282=46449","G282=2725","F283=16907","G283=1000","F285=43406","G285=3456","F286=55099","G286=4256",
"F287=58683","G287=3686","F288=59437","G288=3980","F289=63078","G289=3587","F290=17341","G290=1021",
"F292=83621","G292=7332","F293=83472","G293=5763","F294=93781","G294=4920","F295=16835","G295=1018",
"F299=66523","G299=6856","F300=72795","G300=6943","F301=72489","G301=5700","F302=79075","G302=5650",
"F380=346023","F385=151271","F406=171785","F412=177048","F413=149967","F414=149298","F415=309111",
"F417=344745","G417=593","F421=153632","G421=536","F426=107845","F428=102735","F429=178966","F432=122926",
"G432=434","F433=87464","F434=79496","F435=83681","F436=145556","F439=89495","F440=73930","F441=68549",
"F442=70768","F443=107096","F444=104394","F446=86966","F447=67166","F448=61791","G448=142","F449=60161",
"F453=70389","F454=52734","F455=50878","F457=79149","F458=70951","F460=58627","F462=43802","F464=73747",
"F467=52533","F474=45603","F478=54896","F482=36435","F483=65436","F489=12932","F490=10820","F492=14088",
"F494=6658","F495=10277","F497=7683","G497=10","F498=9062","F499=11223","F500=7555","F502=10369",
"F503=7372","F504=7403","F505=8401","F506=10598","F507=7038","F508=5541","F509=9427","F510=8037",
"F511=6541","G511=20","A512=44406","B512=392581","C512=7288","D512=44","E512=12536","F512=8101",
"G512=16","A513=44407","B513=392647","C513=10239","D513=66","E513=12540","F513=9587","G513=20",
"A514=44408","B514=392704","C514=5864","D514=57","E514=12540","F514=6155","G514=9","A515=44409",
"B515=392710","C515=1205","D515=6","E515=12540","F515=3536","G515=10"
)

foreach ($entry in $updates) {
    $parts = $entry.Split("=")
    $cellRef = $parts[0]
    $newValue = [double]$parts[1]
    $ws.Range($cellRef).Value = $newValue
}

